$d = $word.ActiveDocument

$pairs = @(
    @("2025-10-20 Monday", "2025-10-21 Tuesday"),
    @("83-24=59", "53-34=19"),
    @("94+1=95", "67+1=68"),
    @("14-5=9", "71-0=71"),
    @("15+74=89", "54+36=90"),
    @("45+4=49", "61-37=24"),
    @("3+13=16", "47-28=19"),
    @("21+0=21", "96-85=11"),
    @("62+3=65", "59-9=50"),
    @("52-5=47", "82-81=1"),
    @("8+37=45", "53-0=53"),
    @("38+50=88", "41+40=81"),
    @("80+1=81", "45-9=36"),
    @("80+12=92", "50+44=94"),
    @("28+26=54", "94-4=90"),
    @("31+58=89", "52-34=18"),
    @("49-34=15", "39+55=94"),
    @("81+2=83", "61-0=61"),
    @("50-44=6", "10+72=82"),
    @("60+13=73", "56+26=82"),
    @("92+5=97", "69+26=95"),
    @("57-38=19", "75-1=74"),
    @("34+7=41", "56-29=27"),
    @("40-19=21", "26+6=32"),
    @("99-47=52", "84-69=15"),
    @("19+33=52", "94-36=58"),
    @("93-39=54", "68-17=51"),
    @("78-77=1", "17-12=5"),
    @("60-32=28", "10+83=93"),
    @("17+34=51", "28+67=95"),
    @("38+32=70", "12+75=87"),
    @("92-16=76", "0+31=31"),
    @("90-43=47", "18+33=51"),
    @("41+47=88", "95+4=99"),
    @("14+43=57", "26+0=26"),
    @("77+21=98", "86-21=65"),
    @("7-6=1", "29+38=67"),
    @("67-47=20", "98-56=42"),
    @("10+82=92", "56-25=31"),
    @("85+7=92", "11-9=2"),
    @("88-29=59", "76-45=31"),
    @("25+45=70", "67-54=13"),
    @("40+7=47", "23+1=24"),
    @("5+58=63", "27+17=44"),
    @("42-22=20", "96-84=12"),
    @("19+51=70", "24+30=54"),
    @("78-51=27", "61+35=96"),
    @("69-19=50", "19-6=13"),
    @("34+0=34", "0+11=11"),
    @("39-0=39", "95-15=80"),
    @("95-48=47", "44-16=28"),
    @("47-17=30", "10+2=12"),
    @("88-7=81", "62-6=56"),
    @("68+10=78", "50+37=87"),
    @("10+64=74", "21+8=29"),
    @("66-9=57", "28-25=3"),
    @("65-46=19", "37+62=99"),
    @("80-78=2", "77-49=28"),
    @("36-27=9", "70+5=75"),
    @("45+44=89", "97-22=75"),
    @("6+80=86", "64-47=17"),
    @("58+30=88", "91-74=17"),
    @("96-39=57", "0+22=22"),
    @("0+18=18", "76+4=80"),
    @("28+13=41", "70-4=66"),
    @("85-14=71", "24+29=53"),
    @("34+24=58", "84-78=6"),
    @("3+6=9", "46-17=29"),
    @("26+55=81", "96-33=63"),
    @("30+29=59", "75-47=28"),
    @("56-51=5", "39+8=47"),
    @("20+73=93", "2+28=30"),
    @("22+62=84", "20-4=16"),
    @("88-45=43", "58-24=34"),
    @("45-45=0", "18+12=30"),
    @("51+34=85", "55+14=69"),
    @("3+66=69", "34-20=14"),
    @("3+79=82", "30-16=14"),
    @("32-15=17", "77-13=64"),
    @("97-20=77", "2+55=57"),
    @("24+47=71", "80-37=43"),
    @("14+85=99", "42+42=84"),
    @("32+12=44", "75-7=68"),
    @("99-61=38", "54-35=19"),
    @("32+63=95", "45+12=57"),
    @("77-3=74", "17+61=78"),
    @("63-21=42", "38+40=78"),
    @("97-4=93", "50-21=29"),
    @("16+71=87", "12+67=79"),
    @("73+24=97", "43+54=97"),
    @("86-48=38", "43+21=64"),
    @("59+12=71", "87-17=70"),
    @("60+12=72", "71-59=12"),
    @("88-18=70", "10+60=70"),
    @("38-15=23", "53-0=53"),
    @("60-12=48", "99-59=40"),
    @("76-42=34", "75-54=21"),
    @("66-12=54", "74+0=74"),
    @("58-25=33", "5+36=41"),
    @("86-41=45", "48-1=47"),
    @("31+68=99", "82-18=64")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}
